$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$ws.Range("C4").Value = 660
$ws.Range("C5").Value = 1100
$ws.Range("C6").Value = 1700
